$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 7.846393942920296
$ws.Range("D2").Value = 5.301539641929957
$ws.Range("E2").Value = 12.08062543287042
$ws.Range("F2").Value = 51.94936872646886
$ws.Range("G2").Value = 3.753220566446671
$ws.Range("J2").Value = 10.62410618641181
$ws.Range("K2").Value = 22.00535953010879
$ws.Range("N2").Value = 22.09022614495332
$ws.Range("B3").Value = 7.778037910281712
$ws.Range("D3").Value = 5.302320742851191
$ws.Range("E3").Value = 12.0588024646656
$ws.Range("F3").Value = 51.77238289265132
$ws.Range("G3").Value = 3.757604127648492
$ws.Range("J3").Value = 10.63435391269296
$ws.Range("K3").Value = 21.7379681992359
$ws.Range("N3").Value = 22.14649415479193
$ws.Range("B4").Value = 7.737673882527902
$ws.Range("D4").Value = 5.30357668652882
$ws.Range("E4").Value = 12.04810920929585
$ws.Range("F4").Value = 51.67682606321051
$ws.Range("G4").Value = 3.76043281157938
$ws.Range("J4").Value = 10.64255097331307
$ws.Range("K4").Value = 21.57887887888087
$ws.Range("N4").Value = 22.18301140341903
$ws.Range("B5").Value = 7.721646691400457
$ws.Range("D5").Value = 5.304284450167593
$ws.Range("E5").Value = 12.0444342314277
$ws.Range("F5").Value = 51.64119857604678
$ws.Range("G5").Value = 3.761620155589489
$ws.Range("J5").Value = 10.64636963171163
$ws.Range("K5").Value = 21.51540260356422
$ws.Range("N5").Value = 22.19838765909616
$ws.Range("B6").Value = 7.719011338292184
$ws.Range("D6").Value = 5.304413829629319
$ws.Range("E6").Value = 12.0438652832399
$ws.Range("F6").Value = 51.63548310337324
$ws.Range("G6").Value = 3.761819409044581
$ws.Range("J6").Value = 10.64703258131434
$ws.Range("K6").Value = 21.50494629041052
$ws.Range("N6").Value = 22.20097078840556
$ws.Range("B7").Value = 7.737456005241769
$ws.Range("D7").Value = 5.303585437308387
$ws.Range("E7").Value = 12.04805688099631
$ws.Range("F7").Value = 51.67633214872415
$ws.Range("G7").Value = 3.760448684104177
$ws.Range("J7").Value = 10.64260053754389
$ws.Range("K7").Value = 21.57801723836799
$ws.Range("N7").Value = 22.18321676784039
$ws.Range("B8").Value = 7.822501967700873
$ws.Range("D8").Value = 5.3016482923838
$ws.Range("E8").Value = 12.07254065942579
$ws.Range("F8").Value = 51.88563274813158
$ws.Range("G8").Value = 3.754703636826991
$ws.Range("J8").Value = 10.62724383780338
$ws.Range("K8").Value = 21.91215185079978
$ws.Range("N8").Value = 22.1092189154784
$ws.Range("B9").Value = 8.001225602263949
$ws.Range("D9").Value = 5.303975494583981
$ws.Range("E9").Value = 12.14190808251408
$ws.Range("F9").Value = 52.39926709205334
$ws.Range("G9").Value = 3.744519415204079
$ws.Range("J9").Value = 10.61227264386062
$ws.Range("K9").Value = 22.60439249953157
$ws.Range("N9").Value = 21.97971887317831
$ws.Range("B10").Value = 8.138703665356914
$ws.Range("D10").Value = 5.309371860153774
$ws.Range("E10").Value = 12.20569992693699
$ws.Range("F10").Value = 52.83815285714228
$ws.Range("G10").Value = 3.737687454698942
$ws.Range("J10").Value = 10.61053905229597
$ws.Range("K10").Value = 23.13087411600136
$ws.Range("N10").Value = 21.89407749968791
$ws.Range("B11").Value = 8.202339912286906
$ws.Range("D11").Value = 5.312616015878035
$ws.Range("E11").Value = 12.23745359240382
$ws.Range("F11").Value = 53.05079448607989
$ws.Range("G11").Value = 3.734718693310946
$ws.Range("J11").Value = 10.61176715497323
$ws.Range("K11").Value = 23.37325792008862
$ws.Range("N11").Value = 21.85717796369033
$ws.Range("B12").Value = 8.226574075508813
$ws.Range("D12").Value = 5.313956938535195
$ws.Range("E12").Value = 12.24986573550107
$ws.Range("F12").Value = 53.13314654039295
$ws.Range("G12").Value = 3.733614356718241
$ws.Range("J12").Value = 10.61252232986145
$ws.Range("K12").Value = 23.46537256341244
$ws.Range("N12").Value = 21.84350112922971
$ws.Range("B13").Value = 8.221349082366606
$ws.Range("D13").Value = 5.313663164105209
$ws.Range("E13").Value = 12.24717540754302
$ws.Range("F13").Value = 53.11532984647588
$ws.Range("G13").Value = 3.733851313882477
$ws.Range("J13").Value = 10.61234678593677
$ws.Range("K13").Value = 23.44552074087472
$ws.Range("N13").Value = 21.84643350361039
$ws.Range("B14").Value = 8.204331057559328
$ws.Range("D14").Value = 5.312724087440057
$ws.Range("E14").Value = 12.23846700974578
$ws.Range("F14").Value = 53.05753321139684
$ws.Range("G14").Value = 3.734627441462395
$ws.Range("J14").Value = 10.61182346907395
$ws.Range("K14").Value = 23.38083013765049
$ws.Range("N14").Value = 21.85604682484821
$ws.Range("B15").Value = 8.193924171410412
$ws.Range("D15").Value = 5.312163484970449
$ws.Range("E15").Value = 12.23318318512904
$ws.Range("F15").Value = 53.02236812329119
$ws.Range("G15").Value = 3.735105425440686
$ws.Range("J15").Value = 10.61154070553976
$ws.Range("K15").Value = 23.34124561505334
$ws.Range("N15").Value = 21.86197384210299
$ws.Range("B16").Value = 8.134565144778747
$ws.Range("D16").Value = 5.309175639739943
$ws.Range("E16").Value = 12.20367931703229
$ws.Range("F16").Value = 52.82451439483259
$ws.Range("G16").Value = 3.737884258115799
$ws.Range("J16").Value = 10.61049938767161
$ws.Range("K16").Value = 23.1150844698678
$ws.Range("N16").Value = 21.89653041568191
$ws.Range("B17").Value = 8.098416146464727
$ws.Range("D17").Value = 5.307544112599627
$ws.Range("E17").Value = 12.18627626917701
$ws.Range("F17").Value = 52.70643823342078
$ws.Range("G17").Value = 3.739624518613049
$ws.Range("J17").Value = 10.61037723106682
$ws.Range("K17").Value = 22.9770195747848
$ws.Range("N17").Value = 21.91825722456613
$ws.Range("B18").Value = 8.077728711783346
$ws.Range("D18").Value = 5.306680063783091
$ws.Range("E18").Value = 12.17652423763012
$ws.Range("F18").Value = 52.63974920019425
$ws.Range("G18").Value = 3.740638574399422
$ws.Range("J18").Value = 10.61049679511676
$ws.Range("K18").Value = 22.89788491079437
$ws.Range("N18").Value = 21.93094771890167
$ws.Range("B19").Value = 8.070742897379574
$ws.Range("D19").Value = 5.306400315275601
$ws.Range("E19").Value = 12.17326679307799
$ws.Range("F19").Value = 52.61738099985052
$ws.Range("G19").Value = 3.740984171094016
$ws.Range("J19").Value = 10.61056987427026
$ws.Range("K19").Value = 22.8711414048685
$ws.Range("N19").Value = 21.93527778600107
$ws.Range("B20").Value = 8.102253605856594
$ws.Range("D20").Value = 5.307710102810117
$ws.Range("E20").Value = 12.18810221615754
$ws.Range("F20").Value = 52.71888110292392
$ws.Range("G20").Value = 3.739437909702349
$ws.Range("J20").Value = 10.6103705873649
$ws.Range("K20").Value = 22.99168881355248
$ws.Range("N20").Value = 21.91592430950401
$ws.Range("B21").Value = 8.209326128783257
$ws.Range("D21").Value = 5.312996874488801
$ws.Range("E21").Value = 12.24101440106839
$ws.Range("F21").Value = 53.07446016321287
$ws.Range("G21").Value = 3.734398935840836
$ws.Range("J21").Value = 10.61196930593362
$ws.Range("K21").Value = 23.39982308057626
$ws.Range("N21").Value = 21.85321511898557
$ws.Range("B22").Value = 8.280089977094025
$ws.Range("D22").Value = 5.31710688456004
$ws.Range("E22").Value = 12.27785258034132
$ws.Range("F22").Value = 53.31749379475459
$ws.Range("G22").Value = 3.731221422456642
$ws.Range("J22").Value = 10.61470515239826
$ws.Range("K22").Value = 23.66844206429492
$ws.Range("N22").Value = 21.81395770566842
$ws.Range("B23").Value = 8.242257097457797
$ws.Range("D23").Value = 5.314853749217458
$ws.Range("E23").Value = 12.25798678394567
$ws.Range("F23").Value = 53.18682206891645
$ws.Range("G23").Value = 3.732906775815116
$ws.Range("J23").Value = 10.61309025125156
$ws.Range("K23").Value = 23.52493109218939
$ws.Range("N23").Value = 21.83475208910702
$ws.Range("B24").Value = 8.100518392623407
$ws.Range("D24").Value = 5.307634828254923
$ws.Range("E24").Value = 12.18727591641723
$ws.Range("F24").Value = 52.71325195804916
$ws.Range("G24").Value = 3.739522233353757
$ws.Range("J24").Value = 10.61037299982862
$ws.Range("K24").Value = 22.98505609452384
$ws.Range("N24").Value = 21.91697839913085
$ws.Range("B25").Value = 7.951710099707321
$ws.Range("D25").Value = 5.302694831651428
$ws.Range("E25").Value = 12.12087361019694
$ws.Range("F25").Value = 52.24939561255224
$ws.Range("G25").Value = 3.747159645457362
$ws.Range("J25").Value = 10.61469753668083
$ws.Range("K25").Value = 22.41363831236709
$ws.Range("N25").Value = 22.01308248625694
